# Remove the six oldest years (2004年-2009年), rows 2-7, shifting all remaining
# data (2010年-2020年) up so it now starts at row 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2:7").Delete()

# Append a new row for 2021年 at the bottom (now row 13, right after 2020年 in
# row 12). Copy row 12's formatting first so the new label cell keeps the same
# bold/centered/bordered style used by the other year labels in column A.
$ws.Cells.Item(12, 1).Copy($ws.Cells.Item(13, 1))
$ws.Cells.Item(13, 1).Value = "2021年"

$ws.Cells.Item(13, 2).Value = 188300.45
$ws.Cells.Item(13, 3).Value = 2868.39
$ws.Cells.Item(13, 5).Value = 62420.53
